$d = $word.ActiveDocument

$replacements = @(
    @("162÷9=", "433÷6="),
    @("143÷4=", "985÷6="),
    @("224÷8=", "301÷9="),
    @("484÷8=", "971÷7="),
    @("336÷9=", "196÷3="),
    @("698÷5=", "551÷7="),
    @("917÷3=", "639÷8="),
    @("929÷8=", "652÷2="),
    @("340÷6=", "936÷8="),
    @("855÷3=", "262÷3="),
    @("457÷2=", "410÷9="),
    @("445÷9=", "484÷2="),
    @("509÷4=", "234÷8="),
    @("578÷2=", "662÷3="),
    @("977÷6=", "889÷2="),
    @("613÷5=", "385÷2="),
    @("684÷8=", "206÷6="),
    @("848÷3=", "796÷3="),
    @("345÷7=", "663÷8="),
    @("500÷8=", "376÷7="),
    @("976÷8=", "344÷4="),
    @("891÷8=", "621÷2="),
    @("472÷5=", "866÷2="),
    @("245÷4=", "895÷3="),
    @("339÷6=", "907÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
